# Update "想去人数" (F column) figures for several events on the
# "展览" and "全部类型" worksheets, reflecting refreshed counts from the
# source generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 4603
    3  = 2512
    6  = 61
    7  = 61
    8  = 229
    10 = 183
    12 = 1723
    13 = 316
    14 = 3814
    15 = 31
    16 = 252
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 4603
    3  = 2512
    7  = 61
    8  = 61
    10 = 229
    12 = 183
    16 = 1723
    17 = 316
    18 = 3814
    19 = 31
    20 = 252
}
foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
